# The document embeds two logos (Pearson logo + BTec logo) twice each
# (once in each of the two headers/footers used across sections). The
# edit swaps the "name" attribute values used for the Pearson logo's
# wp:docPr/pic:cNvPr ("image1.png" -> "image2.png") and for the BTec
# logo's wp:docPr/pic:cNvPr ("image2.jpg" -> "image1.jpg"). The "descr"
# attribute (and the underlying relationship/media part) stays the same
# in both cases - only the display "name" changes.
#
# Word's InlineShape object model does not expose the docPr/cNvPr "name"
# attribute (only AlternativeText/Title), so we go through the document's
# flat OOXML package (WordOpenXML) and patch the exact attribute strings,
# then write the updated package back.

$d = $word.ActiveDocument
$xml = $d.WordOpenXML

# -- Pearson logo: wp:docPr id="2" / id="4", and the two pic:cNvPr id="0" --
$xml = $xml.Replace(
    'wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image1.png"',
    'wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image2.png"')

$xml = $xml.Replace(
    'wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="4" name="image1.png"',
    'wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="4" name="image2.png"')

$xml = $xml.Replace(
    'pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image1.png"',
    'pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image2.png"')

# -- BTec logo: wp:docPr id="1" / id="3", and the two pic:cNvPr id="0" --
$xml = $xml.Replace(
    'wp:docPr descr="BTec_Logo-Orange" id="1" name="image2.jpg"',
    'wp:docPr descr="BTec_Logo-Orange" id="1" name="image1.jpg"')

$xml = $xml.Replace(
    'wp:docPr descr="BTec_Logo-Orange" id="3" name="image2.jpg"',
    'wp:docPr descr="BTec_Logo-Orange" id="3" name="image1.jpg"')

$xml = $xml.Replace(
    'pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image2.jpg"',
    'pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image1.jpg"')

$d.WordOpenXML = $xml
